$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D2").Value = "65.702.78"
$ws.Range("E2").Value = "  +0.49%  "
$ws.Range("D3").Value = "3.177.16"
$ws.Range("E3").Value = "  -4.68%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "571.56"
$ws.Range("E5").Value = "  -0.44%  "
$ws.Range("D6").Value = "171.49"
$ws.Range("E6").Value = "  -3.44%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("E8").Value = "  -2.91%  "
$ws.Range("D9").Value = "3.174.61"
$ws.Range("E9").Value = "  -4.65%  "
$ws.Range("E10").Value = "  -2.18%  "
$ws.Range("D11").Value = "6.61"
$ws.Range("E11").Value = "  -3.39%  "
$ws.Range("E12").Value = "  -3.14%  "
$ws.Range("D13").Value = "3.729.76"
$ws.Range("E13").Value = "  -4.60%  "
$ws.Range("E14").Value = "  +0.71%  "
$ws.Range("D15").Value = "27.24"
$ws.Range("E15").Value = "  -3.93%  "
$ws.Range("D16").Value = "65.629.21"
$ws.Range("E16").Value = "  +0.35%  "
$ws.Range("E17").Value = "  -1.97%  "
$ws.Range("D18").Value = "3.180.07"
$ws.Range("E18").Value = "  -4.69%  "
$ws.Range("E19").Value = "  +0.12%  "
$ws.Range("E20").Value = "  -3.30%  "
$ws.Range("D21").Value = "362.17"
$ws.Range("E21").Value = "  +0.06%  "
$ws.Range("D22").Value = "7.29"
$ws.Range("E22").Value = "  -1.74%  "
$ws.Range("D23").Value = "0.999"
$ws.Range("E23").Value = "  -0.04%  "
$ws.Range("D24").Value = "68.77"
$ws.Range("E24").Value = "  -3.19%  "
$ws.Range("E25").Value = "  -3.73%  "
$ws.Range("D26").Value = "3.309.55"
$ws.Range("E26").Value = "  -4.79%  "
$ws.Range("E27").Value = "  -5.85%  "
$ws.Range("D28").Value = "9.89"
$ws.Range("E28").Value = "  +4.04%  "
$ws.Range("E29").Value = "  -0.79%  "
$ws.Range("E31").Value = "  -0.95%  "
$ws.Range("E32").Value = "  -0.08%  "
$ws.Range("E33").Value = "  -3.20%  "
$ws.Range("D34").Value = "22.10"
$ws.Range("E34").Value = "  -3.37%  "
$ws.Range("E35").Value = "  -2.59%  "
$ws.Range("E36").Value = "  -1.09%  "
$ws.Range("D37").Value = "162.32"
$ws.Range("E37").Value = "  +1.48%  "
$ws.Range("E38").Value = "  -1.35%  "
$ws.Range("E39").Value = "  -1.14%  "
$ws.Range("E40").Value = "  +3.37%  "
$ws.Range("D41").Value = "26.50"
$ws.Range("E41").Value = "  -3.02%  "
$ws.Range("E42").Value = "  +1.30%  "
$ws.Range("D43").Value = "2.650.38"
$ws.Range("E43").Value = "  -1.94%  "
$ws.Range("E44").Value = "  -1.24%  "
$ws.Range("E45").Value = "  -1.24%  "
$ws.Range("D46").Value = "39.83"
$ws.Range("E46").Value = "  +0.02%  "
$ws.Range("E47").Value = "  -0.40%  "
$ws.Range("D48").Value = "327.35"
$ws.Range("E48").Value = "  -1.79%  "
$ws.Range("D49").Value = "23.89"
$ws.Range("E49").Value = "  +0.18%  "
$ws.Range("E50").Value = "  -1.30%  "
$ws.Range("E51").Value = "  -0.54%  "
